# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) figures
# for rows 2-51 of the active sheet to match the latest scrape.
#
# For cells whose new text looks like a plain number (e.g. "1.00", "542.60"),
# the Text format ("@") is applied first so Excel keeps the exact literal
# string instead of silently normalising/rounding it into a real number;
# the cell style is then reset to "Normal" so no stray formatting remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.402.69'
$ws.Range("E2").Value = '  -3.00%  '

$ws.Range("D3").Value = '2.271.58'
$ws.Range("E3").Value = '  -6.35%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.36%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.101'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.09%  '

$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.333'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.38%  '

$ws.Range("D14").Value = '2.677.94'
$ws.Range("E14").Value = '  -6.34%  '

$ws.Range("D15").Value = '58.401.37'
$ws.Range("E15").Value = '  -2.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000132'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.40%  '

$ws.Range("D17").Value = '2.272.07'
$ws.Range("E17").Value = '  -6.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '62.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.169'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.37%  '

$ws.Range("E28").Value = '  -2.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.43%  '

$ws.Range("D30").Value = '0.0₃0717'
$ws.Range("E30").Value = '  -7.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.380'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.52%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.20%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '295.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '139.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0943'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0496'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.548'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.84%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0213'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.07%  '

$ws.Range("E49").Value = '  -0.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '

